$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Sheet "QTY Design": append ChR2 / ChR2 QTY rows (rows 24-25)
#    (written first so the new shared strings get indices 44 and 45)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("QTY Design")

$ws3.Range("A24").Value2 = "ChR2"
$ws3.Range("C24").Value2 = 6.13
$ws3.Range("D24").Value2 = 34.9
$ws3.Range("C24:D24").NumberFormat = "0.00"

$ws3.Range("A25").Value2 = "ChR2 QTY"
$ws3.Range("B25").Value2 = 0.235
$ws3.Range("C25").Value2 = 6.13
$ws3.Range("D25").Value2 = 35.3
$ws3.Range("E25").Value2 = 35.53
$ws3.Range("F25").Value2 = 22.22
$ws3.Range("C25:D25").NumberFormat = "0.00"

# ------------------------------------------------------------------
# 2) Sheet "WT-QTY-EXP": insert "CR / monomer" as new row 13 (shifting
#    the old rows 13-14 down to 14-15) and append "CR / dimer" as the
#    new row 16. (Written so the new shared strings get indices 46
#    and 47, in that order.)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("WT-QTY-EXP")

# -- manually shift old rows 13/14 down to 14/15 (bottom-up copy, so
#    nothing is clobbered before it's read), instead of using
#    Rows.Insert() which leaves stray/unused cellXfs behind --
$ws1.Range("A14:D14").Copy() | Out-Null
$ws1.Range("A15:D15").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(15).RowHeight = 50

$ws1.Range("A13:D13").Copy() | Out-Null
$ws1.Range("A14:D14").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(14).RowHeight = 50

# -- fill the new row 13 numbers now but leave the text label for last
#    so the "CR / dimer" string (row 16, added next) is allocated
#    index 46 and "CR / monomer" becomes index 47 --
$ws1.Range("A13:D13").Copy() | Out-Null
$ws1.Range("A13:D13").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(13).RowHeight = 50
$ws1.Range("B13").Value2 = 0.23499999999999999
$ws1.Range("C13").Value2 = 0.42699999999999999
$ws1.Range("D13").Value2 = 0.44700000000000001

# -- append new row 16 ("CR / dimer"); the shift above already moved
#    everything out of the way, so row 16 is a clean append --
$ws1.Range("A15:D15").Copy() | Out-Null
$ws1.Range("A16:D16").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(16).RowHeight = 50
$ws1.Range("A16").Value2 = "CR`ndimer"
$ws1.Range("B16").Value2 = 0.37
$ws1.Range("C16").Value2 = 0.40899999999999997
$ws1.Range("D16").Value2 = 0.51600000000000001

# -- now set row 13's label, last, so it becomes shared-string index 47 --
$ws1.Range("A13").Value2 = "CR`nmonomer"

# ------------------------------------------------------------------
# 3) Selection on "WT-QTY-EXP" moves from G16 to B13
# ------------------------------------------------------------------
$ws1.Range("B13").Select() | Out-Null

# ------------------------------------------------------------------
# 4) Active tab moves from "OPN" to "QTY Design"
# ------------------------------------------------------------------
$ws3.Activate() | Out-Null

Write-Output "edit complete"
